# Applies referee stats refresh for 2025-12-21 run (KHL referees stats 1369).
# Updates per-game cumulative counters on the "Главные" (Main) and "Линейные" (Linear)
# sheets for referees who officiated additional games, and refreshes the
# "as_of_utc" timestamp on every data row of both sheets.

$wb = $excel.ActiveWorkbook

# --- Sheet "Главные" ---
$ws = $wb.Worksheets.Item("Главные")
$updatesMain = @(
    @{ Row = 2; Cells = @{ C = 37; D = 826; E = 370; F = 456; G = 22.32; H = 10; I = 12.32; J = 145; K = 173; L = 10; V = 16; AA = "2025-12-21 04:14:26" } },
    @{ Row = 3; Cells = @{ AA = "2025-12-21 04:14:26" } },
    @{ Row = 4; Cells = @{ AA = "2025-12-21 04:14:26" } },
    @{ Row = 5; Cells = @{ AA = "2025-12-21 04:14:26" } },
    @{ Row = 6; Cells = @{ C = 34; D = 591; E = 282; F = 309; G = 17.38; H = 8.289999999999999; I = 9.09; J = 121; K = 132; AA = "2025-12-21 04:14:26" } },
    @{ Row = 7; Cells = @{ AA = "2025-12-21 04:14:26" } },
    @{ Row = 8; Cells = @{ AA = "2025-12-21 04:14:26" } },
    @{ Row = 9; Cells = @{ C = 34; D = 525; E = 284; F = 241; G = 15.44; H = 8.35; I = 7.09; J = 137; K = 118; L = 2; V = 26; AA = "2025-12-21 04:14:26" } },
    @{ Row = 10; Cells = @{ C = 23; D = 392; E = 200; F = 192; G = 17.04; H = 8.699999999999999; I = 8.35; J = 100; K = 86; AA = "2025-12-21 04:14:26" } },
    @{ Row = 11; Cells = @{ C = 28; D = 652; E = 304; F = 348; G = 23.29; H = 10.86; I = 12.43; J = 137; K = 129; AA = "2025-12-21 04:14:26" } },
    @{ Row = 12; Cells = @{ AA = "2025-12-21 04:14:26" } },
    @{ Row = 13; Cells = @{ AA = "2025-12-21 04:14:26" } },
    @{ Row = 14; Cells = @{ AA = "2025-12-21 04:14:26" } },
    @{ Row = 15; Cells = @{ AA = "2025-12-21 04:14:26" } },
    @{ Row = 16; Cells = @{ AA = "2025-12-21 04:14:26" } },
    @{ Row = 17; Cells = @{ AA = "2025-12-21 04:14:26" } },
    @{ Row = 18; Cells = @{ AA = "2025-12-21 04:14:26" } },
    @{ Row = 19; Cells = @{ AA = "2025-12-21 04:14:26" } },
    @{ Row = 20; Cells = @{ AA = "2025-12-21 04:14:26" } },
    @{ Row = 21; Cells = @{ AA = "2025-12-21 04:14:26" } },
    @{ Row = 22; Cells = @{ AA = "2025-12-21 04:14:26" } },
    @{ Row = 23; Cells = @{ AA = "2025-12-21 04:14:26" } },
    @{ Row = 24; Cells = @{ AA = "2025-12-21 04:14:26" } },
    @{ Row = 25; Cells = @{ AA = "2025-12-21 04:14:26" } },
    @{ Row = 26; Cells = @{ AA = "2025-12-21 04:14:26" } }
)
foreach ($u in $updatesMain) {
    $r = $u.Row
    foreach ($col in $u.Cells.Keys) {
        $ws.Range("$col$r").Value = $u.Cells[$col]
    }
}

# --- Sheet "Линейные" ---
$ws = $wb.Worksheets.Item("Линейные")
$updatesLinear = @(
    @{ Row = 2; Cells = @{ AA = "2025-12-21 04:14:26" } },
    @{ Row = 3; Cells = @{ AA = "2025-12-21 04:14:26" } },
    @{ Row = 4; Cells = @{ AA = "2025-12-21 04:14:26" } },
    @{ Row = 5; Cells = @{ C = 17; D = 236; E = 124; F = 112; G = 13.88; H = 7.29; I = 6.59; J = 62; K = 56; AA = "2025-12-21 04:14:26" } },
    @{ Row = 6; Cells = @{ AA = "2025-12-21 04:14:26" } },
    @{ Row = 7; Cells = @{ AA = "2025-12-21 04:14:26" } },
    @{ Row = 8; Cells = @{ AA = "2025-12-21 04:14:26" } },
    @{ Row = 9; Cells = @{ AA = "2025-12-21 04:14:26" } },
    @{ Row = 10; Cells = @{ AA = "2025-12-21 04:14:26" } },
    @{ Row = 11; Cells = @{ AA = "2025-12-21 04:14:26" } },
    @{ Row = 12; Cells = @{ AA = "2025-12-21 04:14:26" } },
    @{ Row = 13; Cells = @{ AA = "2025-12-21 04:14:26" } },
    @{ Row = 14; Cells = @{ AA = "2025-12-21 04:14:26" } },
    @{ Row = 15; Cells = @{ AA = "2025-12-21 04:14:26" } },
    @{ Row = 16; Cells = @{ AA = "2025-12-21 04:14:26" } },
    @{ Row = 17; Cells = @{ AA = "2025-12-21 04:14:26" } },
    @{ Row = 18; Cells = @{ AA = "2025-12-21 04:14:26" } },
    @{ Row = 19; Cells = @{ AA = "2025-12-21 04:14:26" } },
    @{ Row = 20; Cells = @{ C = 22; D = 355; E = 200; F = 155; G = 16.14; H = 9.09; I = 7.05; J = 85; K = 75; L = 2; V = 14; AA = "2025-12-21 04:14:26" } },
    @{ Row = 21; Cells = @{ AA = "2025-12-21 04:14:26" } },
    @{ Row = 22; Cells = @{ AA = "2025-12-21 04:14:26" } },
    @{ Row = 23; Cells = @{ AA = "2025-12-21 04:14:26" } },
    @{ Row = 24; Cells = @{ C = 36; D = 634; E = 260; F = 374; G = 17.61; H = 7.22; I = 10.39; J = 115; K = 147; AA = "2025-12-21 04:14:26" } },
    @{ Row = 25; Cells = @{ C = 12; D = 187; E = 111; F = 76; G = 15.58; H = 9.25; I = 6.33; J = 53; K = 33; AA = "2025-12-21 04:14:26" } },
    @{ Row = 26; Cells = @{ C = 32; D = 619; E = 283; F = 336; G = 19.34; I = 10.5; J = 114; K = 108; L = 9; V = 14; AA = "2025-12-21 04:14:26" } }
)
foreach ($u in $updatesLinear) {
    $r = $u.Row
    foreach ($col in $u.Cells.Keys) {
        $ws.Range("$col$r").Value = $u.Cells[$col]
    }
}

Write-Output "Updated stats for Главные and Линейные sheets."
